$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Nature comptable",
    "Designation comptable",
    "Centre de coût",
    "Designation centre de coût",
    "Centre de profit",
    "Designation centre de profit",
    "Montant",
    "Type Piece",
    "Nom",
    "Prenom",
    "Matricule",
    "Periode d'effet",
    "Debut periode",
    "Fin periode",
    "N° piece reference",
    "Utilisateur ecriture",
    "Date piece",
    "Date comptable",
    "Date de saisie",
    "Compte contre partie",
    "Designation compte contre partie",
    "N° Ecriture",
    "Commentaire ecriture",
    "N° contre passation",
    "Commentaire contre passation",
    "Devise",
    "Convertion en euros",
    "Date convertion",
    "Taux convertion",
    "Source convertion",
    "Societe",
    "Designation societe",
    "Unité de quantité",
    "Quantité",
    "Taux unité de quantité",
    "Code mouvement",
    "Designation mouvement"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
